# Applies "feat: updated results for abstract" change:
# Adds two new columns G (num_samples) and H (fractional_uncertainty)
# with a header row and per-row bootstrap results for the
# REGIONWISE_BOOTSTRAP thermometry analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("G1").Value = "num_samples"
$ws.Range("H1").Value = "fractional_uncertainty"

# Copy the formatting of an existing header cell (bold, bordered,
# centered/top-aligned) onto the two new header cells so they match
# the look of the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2-19) ---
# num_samples values
$numSamples = @(1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 993, 1000, 988, 1000, 999, 1000, 995, 1000)

# fractional_uncertainty values
$fracUncertainty = @(
    0.05465967198530568,
    0.1517863286558169,
    0.03076134861937979,
    0.1002143141730702,
    0.03284310601598622,
    0.07607541391771377,
    0.03160184766307855,
    0.05038480137631279,
    0.03108241073896207,
    0.05481633040125691,
    0.1057280547432848,
    0.03084771414731377,
    0.1337567131970655,
    0.05682149308062914,
    0.127338086543459,
    0.06317416027103828,
    0.1301786259478235,
    0.05576595618402488
)

for ($i = 0; $i -lt $numSamples.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $numSamples[$i]
    $ws.Cells.Item($row, 8).Value = $fracUncertainty[$i]
}
